$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 28
$ws.Cells.Item(28, 8).Value = 1791.2307
$ws.Cells.Item(28, 9).Value = 909.25
$ws.Cells.Item(28, 11).Value = 909.25
$ws.Cells.Item(28, 13).Value = -424.25

# Row 40
$ws.Cells.Item(40, 8).Value = 7860.5
$ws.Cells.Item(40, 10).Value = 8980
$ws.Cells.Item(40, 12).Value = 8980
$ws.Cells.Item(40, 14).Value = -9330

# Row 42
$ws.Cells.Item(42, 8).Value = 597.44446
$ws.Cells.Item(42, 9).Value = 759.8570999999999
$ws.Cells.Item(42, 11).Value = 2279.5713
$ws.Cells.Item(42, 13).Value = -2049.5713

# Row 43
$ws.Cells.Item(43, 8).Value = 3592.2856
$ws.Cells.Item(43, 9).Value = 3250.5
$ws.Cells.Item(43, 10).Value = 3649.25
$ws.Cells.Item(43, 11).Value = 3250.5
$ws.Cells.Item(43, 12).Value = 3649.25
$ws.Cells.Item(43, 13).Value = -3181.5
$ws.Cells.Item(43, 14).Value = -3787.25

# Row 129
$ws.Cells.Item(129, 8).Value = 11113315
$ws.Cells.Item(129, 9).Value = 1279
$ws.Cells.Item(129, 10).Value = 17546600
$ws.Cells.Item(129, 11).Value = 3837
$ws.Cells.Item(129, 12).Value = 52639800
$ws.Cells.Item(129, 13).Value = 1163
$ws.Cells.Item(129, 14).Value = -52649800

# Row 132
$ws.Cells.Item(132, 8).Value = 217828.7
$ws.Cells.Item(132, 9).Value = 331734
$ws.Cells.Item(132, 10).Value = 6290.2856
$ws.Cells.Item(132, 11).Value = 995202
$ws.Cells.Item(132, 12).Value = 18870.8568
$ws.Cells.Item(132, 13).Value = -992672
$ws.Cells.Item(132, 14).Value = -23930.8568

# Row 135
$ws.Cells.Item(135, 8).Value = 3479.3845
$ws.Cells.Item(135, 9).Value = 1834.2858
$ws.Cells.Item(135, 11).Value = 16508.5722
$ws.Cells.Item(135, 13).Value = -13973.5722

# Row 141
$ws.Cells.Item(141, 8).Value = 5071.143
$ws.Cells.Item(141, 9).Value = 4666.3335
$ws.Cells.Item(141, 10).Value = 7500
$ws.Cells.Item(141, 11).Value = 13999.0005
$ws.Cells.Item(141, 12).Value = 22500
$ws.Cells.Item(141, 13).Value = -8819.000499999998
$ws.Cells.Item(141, 14).Value = -32860


$ws = $wb.Worksheets.Item("ARM")
# Row 92
$ws.Cells.Item(92, 8).Value = 273474
$ws.Cells.Item(92, 10).Value = 273474
$ws.Cells.Item(92, 12).Value = 273474
$ws.Cells.Item(92, 14).Value = -278466

# Row 97
$ws.Cells.Item(97, 8).Value = 744.8261
$ws.Cells.Item(97, 9).Value = 586.8889
$ws.Cells.Item(97, 10).Value = 1313.4
$ws.Cells.Item(97, 11).Value = 586.8889
$ws.Cells.Item(97, 12).Value = 1313.4
$ws.Cells.Item(97, 13).Value = -90.88890000000004
$ws.Cells.Item(97, 14).Value = -2305.4

# Row 112
$ws.Cells.Item(112, 8).Value = 0
$ws.Cells.Item(112, 10).Value = 0
$ws.Cells.Item(112, 12).Value = 0
$ws.Cells.Item(112, 14).ClearContents()

# Row 122
$ws.Cells.Item(122, 8).Value = 2051.3704
$ws.Cells.Item(122, 9).Value = 1144.5625
$ws.Cells.Item(122, 10).Value = 3370.3635
$ws.Cells.Item(122, 11).Value = 3433.6875
$ws.Cells.Item(122, 12).Value = 10111.0905
$ws.Cells.Item(122, 13).Value = -983.6875
$ws.Cells.Item(122, 14).Value = -15011.0905


$ws = $wb.Worksheets.Item("BSM")
# Row 110
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()

# Row 134
$ws.Cells.Item(134, 8).Value = 1838003.8
$ws.Cells.Item(134, 9).Value = 2167642
$ws.Cells.Item(134, 11).Value = 6502926
$ws.Cells.Item(134, 13).Value = -6500391


$ws = $wb.Worksheets.Item("CRP")
# Row 132
$ws.Cells.Item(132, 8).Value = 3401.5642
$ws.Cells.Item(132, 9).Value = 2511.9092
$ws.Cells.Item(132, 11).Value = 7535.7276
$ws.Cells.Item(132, 13).Value = -5005.7276


$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Cells.Item(5, 8).Value = 179533.67
$ws.Cells.Item(5, 9).Value = 928.65216
$ws.Cells.Item(5, 10).Value = 1001116.8
$ws.Cells.Item(5, 11).Value = 2785.95648
$ws.Cells.Item(5, 12).Value = 3003350.4
$ws.Cells.Item(5, 13).Value = -2673.95648
$ws.Cells.Item(5, 14).Value = -3003574.4

# Row 113
$ws.Cells.Item(113, 8).Value = 500645.7
$ws.Cells.Item(113, 9).Value = 495.2857
$ws.Cells.Item(113, 10).Value = 1667663.4
$ws.Cells.Item(113, 11).Value = 1485.8571
$ws.Cells.Item(113, 12).Value = 5002990.199999999
$ws.Cells.Item(113, 13).Value = 684.1428999999998
$ws.Cells.Item(113, 14).Value = -5007330.199999999

# Row 135
$ws.Cells.Item(135, 8).Value = 179533.67
$ws.Cells.Item(135, 9).Value = 928.65216
$ws.Cells.Item(135, 10).Value = 1001116.8
$ws.Cells.Item(135, 11).Value = 8357.86944
$ws.Cells.Item(135, 12).Value = 9010051.200000001
$ws.Cells.Item(135, 13).Value = -5822.86944
$ws.Cells.Item(135, 14).Value = -9015121.200000001


$ws = $wb.Worksheets.Item("GSM")
# Row 111
$ws.Cells.Item(111, 8).Value = 0
$ws.Cells.Item(111, 10).Value = 0
$ws.Cells.Item(111, 12).Value = 0
$ws.Cells.Item(111, 14).ClearContents()


$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Cells.Item(22, 8).Value = 549.5714
$ws.Cells.Item(22, 9).Value = 600
$ws.Cells.Item(22, 11).Value = 600
$ws.Cells.Item(22, 13).Value = -305

# Row 27
$ws.Cells.Item(27, 8).Value = 549.5714
$ws.Cells.Item(27, 9).Value = 600
$ws.Cells.Item(27, 11).Value = 600
$ws.Cells.Item(27, 13).Value = -493

# Row 68
$ws.Cells.Item(68, 8).Value = 2199.7273
$ws.Cells.Item(68, 9).Value = 2199.7273
$ws.Cells.Item(68, 11).Value = 2199.7273
$ws.Cells.Item(68, 13).Value = -1450.7273

# Row 71
$ws.Cells.Item(71, 8).Value = 2199.7273
$ws.Cells.Item(71, 9).Value = 2199.7273
$ws.Cells.Item(71, 11).Value = 10998.6365
$ws.Cells.Item(71, 13).Value = -7254.636500000001

# Row 88
$ws.Cells.Item(88, 8).Value = 23792.334
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 23792.334
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 12).Value = 23792.334
$ws.Cells.Item(88, 13).ClearContents()
$ws.Cells.Item(88, 14).Value = -24648.334

# Row 91
$ws.Cells.Item(91, 8).Value = 23792.334
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 10).Value = 23792.334
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 12).Value = 23792.334
$ws.Cells.Item(91, 13).ClearContents()
$ws.Cells.Item(91, 14).Value = -26756.334

# Row 93
$ws.Cells.Item(93, 8).Value = 1539.381
$ws.Cells.Item(93, 9).Value = 1808.2222
$ws.Cells.Item(93, 10).Value = 1337.75
$ws.Cells.Item(93, 11).Value = 1808.2222
$ws.Cells.Item(93, 12).Value = 1337.75
$ws.Cells.Item(93, 13).Value = -560.2221999999999
$ws.Cells.Item(93, 14).Value = -3833.75

# Row 110
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()


$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Cells.Item(62, 8).Value = 11957.917
$ws.Cells.Item(62, 9).Value = 13435.5
$ws.Cells.Item(62, 10).Value = 11219.125
$ws.Cells.Item(62, 11).Value = 13435.5
$ws.Cells.Item(62, 12).Value = 11219.125
$ws.Cells.Item(62, 13).Value = -12811.5
$ws.Cells.Item(62, 14).Value = -12467.125

# Row 65
$ws.Cells.Item(65, 8).Value = 11957.917
$ws.Cells.Item(65, 9).Value = 13435.5
$ws.Cells.Item(65, 10).Value = 11219.125
$ws.Cells.Item(65, 11).Value = 67177.5
$ws.Cells.Item(65, 12).Value = 56095.625
$ws.Cells.Item(65, 13).Value = -64057.5
$ws.Cells.Item(65, 14).Value = -62335.625

# Row 81
$ws.Cells.Item(81, 8).Value = 797.1667
$ws.Cells.Item(81, 9).Value = 756.6
$ws.Cells.Item(81, 10).Value = 1000
$ws.Cells.Item(81, 11).Value = 1513.2
$ws.Cells.Item(81, 12).Value = 2000
$ws.Cells.Item(81, 13).Value = -452.2
$ws.Cells.Item(81, 14).Value = -4122

# Row 84
$ws.Cells.Item(84, 8).Value = 797.1667
$ws.Cells.Item(84, 9).Value = 756.6
$ws.Cells.Item(84, 10).Value = 1000
$ws.Cells.Item(84, 11).Value = 7566
$ws.Cells.Item(84, 12).Value = 10000
$ws.Cells.Item(84, 13).Value = -2262
$ws.Cells.Item(84, 14).Value = -20608

# Row 88
$ws.Cells.Item(88, 8).Value = 0
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 11).Value = 0
$ws.Cells.Item(88, 13).ClearContents()

# Row 91
$ws.Cells.Item(91, 8).Value = 0
$ws.Cells.Item(91, 9).Value = 0
$ws.Cells.Item(91, 11).Value = 0
$ws.Cells.Item(91, 13).ClearContents()

# Row 96
$ws.Cells.Item(96, 8).Value = 2061.4
$ws.Cells.Item(96, 9).Value = 2071.5
$ws.Cells.Item(96, 10).Value = 2059.8462
$ws.Cells.Item(96, 11).Value = 2071.5
$ws.Cells.Item(96, 12).Value = 2059.8462
$ws.Cells.Item(96, 13).Value = -698.5
$ws.Cells.Item(96, 14).Value = -4805.8462

# Row 107
$ws.Cells.Item(107, 8).Value = 605.8125
$ws.Cells.Item(107, 9).Value = 401
$ws.Cells.Item(107, 10).Value = 947.1667
$ws.Cells.Item(107, 11).Value = 1203
$ws.Cells.Item(107, 12).Value = 2841.5001
$ws.Cells.Item(107, 13).Value = 717
$ws.Cells.Item(107, 14).Value = -6681.5001

